# Auto-generated Excel COM-interop script applying the Durandal_Profits diff
# Updates market price / profit columns (H-N) for specific Leve rows across sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 934.7143
$ws.Range("I41").Value = 986.5454999999999
$ws.Range("J41").Value = 744.6667
$ws.Range("K41").Value = 986.5454999999999
$ws.Range("L41").Value = 744.6667
$ws.Range("M41").Value = -546.5454999999999
$ws.Range("N41").Value = -1624.6667

$ws.Range("H42").Value = 242.81818
$ws.Range("I42").Value = 278.2
$ws.Range("J42").Value = 213.33333
$ws.Range("K42").Value = 834.5999999999999
$ws.Range("L42").Value = 639.99999
$ws.Range("M42").Value = -604.5999999999999
$ws.Range("N42").Value = -1099.99999

$ws.Range("H88").Value = 4903460.5
$ws.Range("I88").Value = 1800
$ws.Range("J88").Value = 29411764
$ws.Range("K88").Value = 1800
$ws.Range("L88").Value = 29411764
$ws.Range("M88").Value = -1394
$ws.Range("N88").Value = -29412576

$ws.Range("H91").Value = 4903460.5
$ws.Range("I91").Value = 1800
$ws.Range("J91").Value = 29411764
$ws.Range("K91").Value = 1800
$ws.Range("L91").Value = 29411764
$ws.Range("M91").Value = -396
$ws.Range("N91").Value = -29414572

$ws.Range("H98").Value = 3972727.5
$ws.Range("I98").Value = 5042.0435
$ws.Range("J98").Value = 22224080
$ws.Range("K98").Value = 5042.0435
$ws.Range("L98").Value = 22224080
$ws.Range("M98").Value = -3544.0435
$ws.Range("N98").Value = -22227076

$ws.Range("H122").Value = 3972727.5
$ws.Range("I122").Value = 5042.0435
$ws.Range("J122").Value = 22224080
$ws.Range("K122").Value = 15126.1305
$ws.Range("L122").Value = 66672240
$ws.Range("M122").Value = -12676.1305
$ws.Range("N122").Value = -66677140

$ws.Range("H138").Value = 3166.8289
$ws.Range("I138").Value = 2161.5293
$ws.Range("J138").Value = 3980.6428
$ws.Range("K138").Value = 6484.5879
$ws.Range("L138").Value = 11941.9284
$ws.Range("M138").Value = -1344.5879
$ws.Range("N138").Value = -22221.9284

$ws.Range("H139").Value = 77600
$ws.Range("J139").Value = 77600
$ws.Range("L139").Value = 77600
$ws.Range("N139").Value = -87880

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 496533.12
$ws.Range("I32").Value = 7092.3086
$ws.Range("J32").Value = 5452121.5
$ws.Range("K32").Value = 7092.3086
$ws.Range("L32").Value = 5452121.5
$ws.Range("M32").Value = -6805.3086
$ws.Range("N32").Value = -5452695.5

$ws.Range("H37").Value = 5191.3335
$ws.Range("J37").Value = 8038
$ws.Range("L37").Value = 8038
$ws.Range("N37").Value = -8584

$ws.Range("H138").Value = 69157.14
$ws.Range("J138").Value = 69157.14
$ws.Range("L138").Value = 69157.14
$ws.Range("N138").Value = -79437.14

$ws.Range("H139").Value = 54922.5
$ws.Range("J139").Value = 54922.5
$ws.Range("L139").Value = 54922.5
$ws.Range("N139").Value = -65202.5

$ws.Range("H140").Value = 102832.9
$ws.Range("J140").Value = 102832.9
$ws.Range("L140").Value = 102832.9
$ws.Range("N140").Value = -113192.9

$ws.Range("H141").Value = 64766.668
$ws.Range("J141").Value = 64766.668
$ws.Range("L141").Value = 64766.668
$ws.Range("N141").Value = -75126.66800000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1544.9524
$ws.Range("I94").Value = 1455.5
$ws.Range("J94").Value = 2081.6667
$ws.Range("K94").Value = 1455.5
$ws.Range("L94").Value = 2081.6667
$ws.Range("M94").Value = -1004.5
$ws.Range("N94").Value = -2983.6667

$ws.Range("H105").Value = 1697.1875
$ws.Range("I105").Value = 1696.6666
$ws.Range("J105").Value = 1700
$ws.Range("K105").Value = 1696.6666
$ws.Range("L105").Value = 1700
$ws.Range("M105").Value = 50.33339999999998
$ws.Range("N105").Value = -5194

$ws.Range("H107").Value = 3951.394
$ws.Range("I107").Value = 4481.2593
$ws.Range("J107").Value = 1567
$ws.Range("K107").Value = 4481.2593
$ws.Range("L107").Value = 1567
$ws.Range("M107").Value = -2561.2593
$ws.Range("N107").Value = -5407

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9974.078
$ws.Range("I31").Value = 3494.606
$ws.Range("J31").Value = 14833.682
$ws.Range("K31").Value = 3494.606
$ws.Range("L31").Value = 14833.682
$ws.Range("M31").Value = -3199.606
$ws.Range("N31").Value = -15423.682

$ws.Range("H34").Value = 9974.078
$ws.Range("I34").Value = 3494.606
$ws.Range("J34").Value = 14833.682
$ws.Range("K34").Value = 3494.606
$ws.Range("L34").Value = 14833.682
$ws.Range("M34").Value = -3292.606
$ws.Range("N34").Value = -15237.682

$ws.Range("H51").Value = 9098
$ws.Range("J51").Value = 9098
$ws.Range("L51").Value = 9098
$ws.Range("N51").Value = -10570

$ws.Range("H60").Value = 6250.5
$ws.Range("J60").Value = 8009
$ws.Range("L60").Value = 8009
$ws.Range("N60").Value = -9031

$ws.Range("H61").Value = 9098
$ws.Range("J61").Value = 9098
$ws.Range("L61").Value = 9098
$ws.Range("N61").Value = -9794

$ws.Range("H74").Value = 18227.4
$ws.Range("J74").Value = 18227.4
$ws.Range("L74").Value = 18227.4
$ws.Range("N74").Value = -19975.4

$ws.Range("H77").Value = 18227.4
$ws.Range("J77").Value = 18227.4
$ws.Range("L77").Value = 54682.2
$ws.Range("N77").Value = -63418.2

$ws.Range("H99").Value = 1834.9429
$ws.Range("I99").Value = 1907.4667
$ws.Range("J99").Value = 1399.8
$ws.Range("K99").Value = 1907.4667
$ws.Range("L99").Value = 1399.8
$ws.Range("M99").Value = -409.4666999999999
$ws.Range("N99").Value = -4395.8

$ws.Range("H105").Value = 1001.3571
$ws.Range("I105").Value = 963
$ws.Range("J105").Value = 1500
$ws.Range("K105").Value = 963
$ws.Range("L105").Value = 1500
$ws.Range("M105").Value = 784
$ws.Range("N105").Value = -4994

$ws.Range("H126").Value = 1834.9429
$ws.Range("I126").Value = 1907.4667
$ws.Range("J126").Value = 1399.8
$ws.Range("K126").Value = 5722.4001
$ws.Range("L126").Value = 4199.4
$ws.Range("M126").Value = -3252.4001
$ws.Range("N126").Value = -9139.4

$ws.Range("H132").Value = 1086.9722
$ws.Range("I132").Value = 673.5172
$ws.Range("K132").Value = 2020.5516
$ws.Range("M132").Value = 509.4484

$ws.Range("H135").Value = 50833.332
$ws.Range("J135").Value = 50833.332
$ws.Range("L135").Value = 50833.332
$ws.Range("N135").Value = -60973.332

$ws.Range("H138").Value = 48144.445
$ws.Range("J138").Value = 48144.445
$ws.Range("L138").Value = 48144.445
$ws.Range("N138").Value = -58424.445

$ws.Range("H140").Value = 89300
$ws.Range("J140").Value = 89300
$ws.Range("L140").Value = 89300
$ws.Range("N140").Value = -99660

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 363.86957
$ws.Range("I23").Value = 197.64285
$ws.Range("J23").Value = 622.44446
$ws.Range("K23").Value = 592.9285500000001
$ws.Range("L23").Value = 1867.33338
$ws.Range("M23").Value = -357.9285500000001
$ws.Range("N23").Value = -2337.33338

$ws.Range("H121").Value = 37017.656
$ws.Range("I121").Value = 317.4
$ws.Range("J121").Value = 40687.68
$ws.Range("K121").Value = 952.1999999999999
$ws.Range("L121").Value = 122063.04
$ws.Range("M121").Value = 357.8000000000001
$ws.Range("N121").Value = -124683.04

$ws.Range("H122").Value = 456
$ws.Range("I122").Value = 346.61765
$ws.Range("J122").Value = 1199.8
$ws.Range("K122").Value = 3119.55885
$ws.Range("L122").Value = 10798.2
$ws.Range("M122").Value = -669.5588500000003
$ws.Range("N122").Value = -15698.2

$ws.Range("H129").Value = 17544968
$ws.Range("I129").Value = 1060.3334
$ws.Range("J129").Value = 25642156
$ws.Range("K129").Value = 3181.0002
$ws.Range("L129").Value = 76926468
$ws.Range("M129").Value = 1818.9998
$ws.Range("N129").Value = -76936468

$ws.Range("H134").Value = 2946.7441
$ws.Range("I134").Value = 2136.2334
$ws.Range("J134").Value = 4817.154
$ws.Range("K134").Value = 6408.7002
$ws.Range("L134").Value = 14451.462
$ws.Range("M134").Value = -1338.7002
$ws.Range("N134").Value = -24591.462

$ws.Range("H139").Value = 2348.9111
$ws.Range("I139").Value = 1217.5
$ws.Range("K139").Value = 3652.5
$ws.Range("M139").Value = 1487.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H138").Value = 69233.336
$ws.Range("J138").Value = 69233.336
$ws.Range("L138").Value = 69233.336
$ws.Range("N138").Value = -79513.336

$ws.Range("H139").Value = 52423.168
$ws.Range("J139").Value = 52423.168
$ws.Range("L139").Value = 52423.168
$ws.Range("N139").Value = -62703.168

$ws.Range("H140").Value = 89899
$ws.Range("J140").Value = 89899
$ws.Range("L140").Value = 89899
$ws.Range("N140").Value = -100259

$ws.Range("H141").Value = 48000
$ws.Range("J141").Value = 48000
$ws.Range("L141").Value = 48000
$ws.Range("N141").Value = -58360

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 710
$ws.Range("I22").Value = 509.0909
$ws.Range("J22").Value = 955.55554
$ws.Range("K22").Value = 509.0909
$ws.Range("L22").Value = 955.55554
$ws.Range("M22").Value = -214.0909
$ws.Range("N22").Value = -1545.55554

$ws.Range("H27").Value = 710
$ws.Range("I27").Value = 509.0909
$ws.Range("J27").Value = 955.55554
$ws.Range("K27").Value = 509.0909
$ws.Range("L27").Value = 955.55554
$ws.Range("M27").Value = -402.0909
$ws.Range("N27").Value = -1169.55554

$ws.Range("H46").Value = 687.75
$ws.Range("I46").Value = 683.6667
$ws.Range("K46").Value = 683.6667
$ws.Range("M46").Value = -495.6667

$ws.Range("H139").Value = 64540
$ws.Range("J139").Value = 79425
$ws.Range("L139").Value = 79425
$ws.Range("N139").Value = -89705

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 984.2308
$ws.Range("I136").Value = 618.1905
$ws.Range("J136").Value = 1411.2778
$ws.Range("K136").Value = 1854.5715
$ws.Range("L136").Value = 4233.8334
$ws.Range("M136").Value = 695.4285
$ws.Range("N136").Value = -9333.8334

$ws.Range("H138").Value = 78525
$ws.Range("J138").Value = 93033.336
$ws.Range("L138").Value = 93033.336
$ws.Range("N138").Value = -103313.336

$ws.Range("H139").Value = 53376.875
$ws.Range("J139").Value = 53376.875
$ws.Range("L139").Value = 53376.875
$ws.Range("N139").Value = -63656.875

$ws.Range("H141").Value = 78653.44500000001
$ws.Range("J141").Value = 78653.44500000001
$ws.Range("L141").Value = 78653.44500000001
$ws.Range("N141").Value = -89013.44500000001

Write-Output "Applied all updates"
